$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the current row 25 (the "ERROR_TEST" row),
# which pushes that row down to row 26 and makes room for the new
# "ERROR_LOGOUT_FAILED" entry at row 25.
$ws.Rows.Item(25).Insert()

$ws.Cells.Item(25, 1).Value = 1023
$ws.Cells.Item(25, 2).Value = "ERROR_LOGOUT_FAILED"
$ws.Cells.Item(25, 3).Value = "general"
$ws.Cells.Item(25, 4).Value = "注销失败"

$ws.Range("E26").Select()
